$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: double the unit prices in D29 and D30
$ws.Range("D29").Value = 112.4
$ws.Range("D30").Value = 187.2
